$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Wolverine ACA Merrell LS" (shared by cell E5) is being renamed to "Abbott"
$ws.Range("E5").Value = "Abbott"

# New value "DTC" added to cell E6 (previously empty)
$ws.Range("E6").Value = "DTC"

# Move/restore the active selection to E10
$ws.Range("E10").Select() | Out-Null
